$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.139.84"
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("D3").Value = "2.613.22"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'568.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.63%  "
$ws.Range("D6").Value = "'145.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "'0.599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.97%  "
$ws.Range("D9").Value = "2.624.90"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").Value = "'6.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("E11").Value = "  +3.51%  "
$ws.Range("E12").Value = "  +9.57%  "
$ws.Range("D13").Value = "'0.343"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("D14").Value = "3.077.69"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").Value = "60.086.44"
$ws.Range("E15").Value = "  +2.71%  "
$ws.Range("D16").Value = "'22.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.36%  "
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").Value = "2.639.47"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D19").Value = "'4.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "'340.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("D21").Value = "'10.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'6.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.62%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'65.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").Value = "'0.448"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.78%  "
$ws.Range("E26").Value = "  +3.04%  "
$ws.Range("D27").Value = "'0.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'7.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.65%  "
$ws.Range("D29").Value = "0.0₃0791"
$ws.Range("E29").Value = "  +7.22%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +3.79%  "
$ws.Range("D32").Value = "'6.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("D33").Value = "'159.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("D34").Value = "'19.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.34%  "
$ws.Range("D35").Value = "'4.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.43%  "
$ws.Range("D36").Value = "'1.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.82%  "
$ws.Range("D37").Value = "'0.885"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.77%  "
$ws.Range("D38").Value = "'0.879"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.76%  "
$ws.Range("D39").Value = "'37.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("D40").Value = "'1.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.99%  "
$ws.Range("D41").Value = "'296.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.42%  "
$ws.Range("D42").Value = "'3.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("D43").Value = "'0.994"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").Value = "'0.0978"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.599"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "'10.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'19.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.51%  "
$ws.Range("D49").Value = "'126.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.30%  "
$ws.Range("E50").Value = "  +3.33%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'4.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.32%  "
